$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append 5 new rows (85-89) to the "Dominios SDS" table for the new
# "Ruptura Vitrine" solution, mirroring the layout of the existing rows ---

$newRows = @(
    @("Ruptura Vitrine", "stg_ruptura.dtsx", "LOJACORP", "dbo. spETL_source_Sku",      "MIS_DW", "dbo.stg_corp_sku",      "stg_sku_ruptura_vitrine.dtsx"),
    @("Ruptura Vitrine", "stg_ruptura.dtsx", "LOJACORP", "dbo.spETL_source_Produto",   "MIS_DW", "dbo.stg_corp_product",  "stg_produto_ruptura_vitrine.dtsx"),
    @("Ruptura Vitrine", "stg_ruptura.dtsx", "LOJACORP", "dbo.spETL_source_Categoria", "MIS_DW", "dbo.stg_corp_categoria","stg_categoria_ruptura_vitrine.dtsx"),
    @("Ruptura Vitrine", "stg_ruptura.dtsx", "LOJACORP", "dbo.spETL_source_SkuKit",    "MIS_DW", "stg_corp_kit",          "stg_kit_ruptura_vitrine.dtsx"),
    @("Ruptura Vitrine", "stg_ruptura.dtsx", "LOJACORP", "dbo.Produto",                "MIS_DW", "stg_loja_product",      "stg_product_ruptura_vitrine.dtsx")
)

$startRow = 85
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# --- Grow the "Tabela1" table (and its AutoFilter) so it covers the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G89"))

# --- Update the hidden _FilterDatabase defined name to match the new extent ---
$wb.Names.Item(1).RefersTo = "='Dominios SDS '!`$A`$1:`$F`$89"

# --- Column G widened (bestFit) to accommodate the new, longer DTSX names ---
$ws.Columns.Item(7).ColumnWidth = 31.3

# --- Restore the active-cell selection recorded for the bottom-right pane ---
$ws.Range("D78").Select() | Out-Null
